# Add missing NIFES columns (HCB + DDT)
# The underlying data fix: the "for R" sheet had a stray/duplicate row
# (HG / Fisk, lever / Gadus morhua / LI / EQS=20) that needs to be removed
# so the NIFES (HCB/DDT) rows line up correctly below. Deleting that row
# shifts every row below it up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("for R")
$ws.Activate()

# Delete entire row 4 - shifts rows 5:34 up to 4:33
$ws.Rows.Item(4).Delete()

# Restore the selection that was left after the edit
$ws.Range("A6").Select()
